# Applies the cryptos-list price/volume/coin-order refresh described by the commit
# "Updated cryptos list on Sun May 26 02:41:21 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D sometimes hold numeric-looking text (e.g. "68.900.18" or
# "600.46") that must stay plain text instead of being auto-converted to a
# number by Excel. Force text via NumberFormat, assign, then restore the
# default "Normal" style so no visible formatting change is left behind.
function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "68.900.18"
$ws.Range("E2").Value = "  +0.26%  "
# Row 3
$ws.Range("D3").Value = "3.739.95"
$ws.Range("E3").Value = "  +0.45%  "
# Row 4
$ws.Range("E4").Value = "  +0.08%  "
# Row 5
Set-TextValue "D5" "600.46"
$ws.Range("E5").Value = "  +0.03%  "
# Row 6
Set-TextValue "D6" "165.33"
$ws.Range("E6").Value = "  -2.14%  "
# Row 7
$ws.Range("D7").Value = "3.739.90"
$ws.Range("E7").Value = "  +0.52%  "
# Row 8
$ws.Range("E8").Value = "  -0.02%  "
# Row 9
Set-TextValue "D9" "0.539"
$ws.Range("E9").Value = "  +0.94%  "
# Row 10
Set-TextValue "D10" "0.170"
$ws.Range("E10").Value = "  +4.52%  "
# Row 11
Set-TextValue "D11" "6.42"
$ws.Range("E11").Value = "  +1.17%  "
# Row 12
Set-TextValue "D12" "0.459"
$ws.Range("E12").Value = "  -0.44%  "
# Row 13
Set-TextValue "D13" "37.68"
$ws.Range("E13").Value = "  -0.97%  "
# Row 14
Set-TextValue "D14" "0.0000247"
$ws.Range("E14").Value = "  +1.12%  "
# Row 15
$ws.Range("D15").Value = "4.370.12"
$ws.Range("E15").Value = "  +0.45%  "
# Row 16
$ws.Range("D16").Value = "3.743.76"
$ws.Range("E16").Value = "  +0.41%  "
# Row 17
$ws.Range("D17").Value = "69.054.65"
$ws.Range("E17").Value = "  +0.48%  "
# Row 18
Set-TextValue "D18" "7.43"
$ws.Range("E18").Value = "  +2.22%  "
# Row 19
Set-TextValue "D19" "17.49"
$ws.Range("E19").Value = "  +1.97%  "
# Row 20
Set-TextValue "D20" "0.113"
$ws.Range("E20").Value = "  -1.99%  "
# Row 21
Set-TextValue "D21" "11.11"
$ws.Range("E21").Value = "  +5.38%  "
# Row 22
Set-TextValue "D22" "491.77"
$ws.Range("E22").Value = "  -0.40%  "
# Row 23
Set-TextValue "D23" "0.726"
$ws.Range("E23").Value = "  +0.04%  "
# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D24" "84.53"
$ws.Range("E24").Value = "  -0.66%  "
# Row 25
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D25" "0.0000147"
$ws.Range("E25").Value = "  +2.94%  "
# Row 26
Set-TextValue "D26" "2.26"
$ws.Range("E26").Value = "  -1.98%  "
# Row 27
Set-TextValue "D27" "12.31"
$ws.Range("E27").Value = "  -1.19%  "
# Row 28
Set-TextValue "D28" "10.02"
$ws.Range("E28").Value = "  -0.97%  "
# Row 29
$ws.Range("E29").Value = "  -0.01%  "
# Row 30
Set-TextValue "D30" "2.96"
$ws.Range("E30").Value = "  +0.03%  "
# Row 31
Set-TextValue "D31" "8.19"
$ws.Range("E31").Value = "  +3.26%  "
# Row 32
Set-TextValue "D32" "2.43"
$ws.Range("E32").Value = "  -6.82%  "
# Row 33
$ws.Range("B33").Value = "WrappedeETH"
$ws.Range("C33").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D33").Value = "3.890.61"
$ws.Range("E33").Value = "  +0.52%  "
# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D34" "31.48"
$ws.Range("E34").Value = "  +0.08%  "
# Row 35
$ws.Range("D35").Value = "3.680.11"
$ws.Range("E35").Value = "  +0.58%  "
# Row 36
$ws.Range("E36").Value = "  -0.53%  "
# Row 37
Set-TextValue "D37" "5.91"
$ws.Range("E37").Value = "  +1.37%  "
# Row 38
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D38" "1.01"
$ws.Range("E38").Value = "  -0.45%  "
# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D39" "0.138"
$ws.Range("E39").Value = "  +5.24%  "
# Row 40
Set-TextValue "D40" "3.16"
$ws.Range("E40").Value = "  +9.68%  "
# Row 41
$ws.Range("E41").Value = "  -0.02%  "
# Row 42
Set-TextValue "D42" "0.323"
$ws.Range("E42").Value = "  -0.45%  "
# Row 43
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D43" "48.56"
$ws.Range("E43").Value = "  -0.50%  "
# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D44" "1.99"
$ws.Range("E44").Value = "  -0.02%  "
# Row 45
Set-TextValue "D45" "424.07"
$ws.Range("E45").Value = "  -3.07%  "
# Row 46
Set-TextValue "D46" "8.42"
# Row 47
$ws.Range("E47").Value = "  +0.00%  "
# Row 48
Set-TextValue "D48" "39.80"
$ws.Range("E48").Value = "  -1.95%  "
# Row 49
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D49" "141.01"
$ws.Range("E49").Value = "  -0.14%  "
# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.775.70"
$ws.Range("E50").Value = "  +0.33%  "
# Row 51
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D51" "1.29"
$ws.Range("E51").Value = "  +5.78%  "
